$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from existing header cell (F1) onto the new header cells
$ws.Range("F1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Set the new header labels
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Set the new metric values
$ws.Range("G2").Value = 0.1228586025167412
$ws.Range("H2").Value = 0.991
